# Update leve-profit figures (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets,
# matching the scheduled-runner refresh from the market-board data source.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 411.9
$ws.Range("I19").Value = 547
$ws.Range("J19").Value = 396.8889
$ws.Range("K19").Value = 547
$ws.Range("L19").Value = 396.8889
$ws.Range("M19").Value = -372
$ws.Range("N19").Value = -746.8888999999999
# Row 138
$ws.Range("H138").Value = 2325.1216
$ws.Range("J138").Value = 3484.15
$ws.Range("L138").Value = 10452.45
$ws.Range("N138").Value = -20732.45

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 39143.69
$ws.Range("J2").Value = 143881.72
$ws.Range("L2").Value = 143881.72
$ws.Range("N2").Value = -144107.72
# Row 32
$ws.Range("H32").Value = 6493.143
$ws.Range("I32").Value = 6429.382
$ws.Range("K32").Value = 6429.382
$ws.Range("M32").Value = -6142.382
# Row 74
$ws.Range("H74").Value = 3190.2144
$ws.Range("I74").Value = 2014.1666
$ws.Range("J74").Value = 4072.25
$ws.Range("K74").Value = 2014.1666
$ws.Range("L74").Value = 4072.25
$ws.Range("M74").Value = -1140.1666
$ws.Range("N74").Value = -5820.25
# Row 77
$ws.Range("H77").Value = 3190.2144
$ws.Range("I77").Value = 2014.1666
$ws.Range("J77").Value = 4072.25
$ws.Range("K77").Value = 10070.833
$ws.Range("L77").Value = 20361.25
$ws.Range("M77").Value = -5702.833000000001
$ws.Range("N77").Value = -29097.25
# Row 110
$ws.Range("H110").Value = 1096.5172
$ws.Range("I110").Value = 819.0417
$ws.Range("K110").Value = 819.0417
$ws.Range("M110").Value = 1225.9583
# Row 116
$ws.Range("H116").Value = 39143.69
$ws.Range("J116").Value = 143881.72
$ws.Range("L116").Value = 143881.72
$ws.Range("N116").Value = -148469.72
# Row 132
$ws.Range("H132").Value = 3015.6572
$ws.Range("I132").Value = 2471.318
$ws.Range("J132").Value = 3936.8462
$ws.Range("K132").Value = 7413.954000000001
$ws.Range("L132").Value = 11810.5386
$ws.Range("M132").Value = -4883.954000000001
$ws.Range("N132").Value = -16870.5386

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 39143.69
$ws.Range("J3").Value = 143881.72
$ws.Range("L3").Value = 143881.72
$ws.Range("N3").Value = -144109.72
# Row 86
$ws.Range("H86").Value = 6611.9443
$ws.Range("J86").Value = 2291.6667
$ws.Range("L86").Value = 2291.6667
$ws.Range("N86").Value = -4537.6667
# Row 89
$ws.Range("H89").Value = 6611.9443
$ws.Range("J89").Value = 2291.6667
$ws.Range("L89").Value = 11458.3335
$ws.Range("N89").Value = -22690.3335
# Row 99
$ws.Range("H99").Value = 19323.227
$ws.Range("J99").Value = 5682.625
$ws.Range("L99").Value = 5682.625
$ws.Range("N99").Value = -8678.625
# Row 134
$ws.Range("H134").Value = 7425.4
$ws.Range("I134").Value = 7769.4736
$ws.Range("K134").Value = 23308.4208
$ws.Range("M134").Value = -20773.4208

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9827.65
$ws.Range("I31").Value = 12146.923
$ws.Range("K31").Value = 12146.923
$ws.Range("M31").Value = -11851.923
# Row 34
$ws.Range("H34").Value = 9827.65
$ws.Range("I34").Value = 12146.923
$ws.Range("K34").Value = 12146.923
$ws.Range("M34").Value = -11944.923
# Row 88
$ws.Range("H88").Value = 40164
$ws.Range("J88").Value = 38196.8
$ws.Range("L88").Value = 38196.8
$ws.Range("N88").Value = -39008.8
# Row 91
$ws.Range("H91").Value = 40164
$ws.Range("J91").Value = 38196.8
$ws.Range("L91").Value = 38196.8
$ws.Range("N91").Value = -41004.8
# Row 99
$ws.Range("H99").Value = 21478998
$ws.Range("I99").Value = 29531276
$ws.Range("J99").Value = 6256.3335
$ws.Range("K99").Value = 29531276
$ws.Range("L99").Value = 6256.3335
$ws.Range("M99").Value = -29529778
$ws.Range("N99").Value = -9252.333500000001
# Row 126
$ws.Range("H126").Value = 21478998
$ws.Range("I126").Value = 29531276
$ws.Range("J126").Value = 6256.3335
$ws.Range("K126").Value = 88593828
$ws.Range("L126").Value = 18769.0005
$ws.Range("M126").Value = -88591358
$ws.Range("N126").Value = -23709.0005
# Row 134
$ws.Range("H134").Value = 5407
$ws.Range("I134").Value = 13925.5
$ws.Range("J134").Value = 1999.6
$ws.Range("K134").Value = 41776.5
$ws.Range("L134").Value = 5998.799999999999
$ws.Range("M134").Value = -39241.5
$ws.Range("N134").Value = -11068.8

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 73.22221999999999
$ws.Range("I2").Value = 89.57143000000001
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 537.42858
$ws.Range("L2").Value = 96
$ws.Range("M2").Value = -424.42858
$ws.Range("N2").Value = -322

$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 17498
$ws.Range("I26").Value = 18246.25
$ws.Range("J26").Value = 16749.75
$ws.Range("K26").Value = 18246.25
$ws.Range("L26").Value = 16749.75
$ws.Range("M26").Value = -17966.25
$ws.Range("N26").Value = -17309.75
# Row 50
$ws.Range("H50").Value = 17498
$ws.Range("I50").Value = 18246.25
$ws.Range("J50").Value = 16749.75
$ws.Range("K50").Value = 18246.25
$ws.Range("L50").Value = 16749.75
$ws.Range("M50").Value = -17748.25
$ws.Range("N50").Value = -17745.75
# Row 58
$ws.Range("H58").Value = 15000
$ws.Range("J58").Value = 15000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15554
# Row 80
$ws.Range("H80").Value = 22376.2
$ws.Range("I80").Value = 22376.2
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 22376.2
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -21378.2
# Row 83
$ws.Range("H83").Value = 22376.2
$ws.Range("I83").Value = 22376.2
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 111881
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -106889
# Row 102
$ws.Range("H102").Value = 8649.333000000001
$ws.Range("I102").Value = 9821.375
$ws.Range("K102").Value = 9821.375
$ws.Range("M102").Value = -8199.375
# Row 113
$ws.Range("H113").Value = 8604.764999999999
$ws.Range("I113").Value = 10548.417
$ws.Range("K113").Value = 10548.417
$ws.Range("M113").Value = -8378.416999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 25748.125
$ws.Range("I22").Value = 40359.6
$ws.Range("J22").Value = 1395.6666
$ws.Range("K22").Value = 40359.6
$ws.Range("L22").Value = 1395.6666
$ws.Range("M22").Value = -40064.6
$ws.Range("N22").Value = -1985.6666
# Row 27
$ws.Range("H27").Value = 25748.125
$ws.Range("I27").Value = 40359.6
$ws.Range("J27").Value = 1395.6666
$ws.Range("K27").Value = 40359.6
$ws.Range("L27").Value = 1395.6666
$ws.Range("M27").Value = -40252.6
$ws.Range("N27").Value = -1609.6666
# Row 40
$ws.Range("H40").Value = 27875
$ws.Range("I40").Value = 29392
$ws.Range("K40").Value = 29392
$ws.Range("M40").Value = -29256
# Row 42
$ws.Range("H42").Value = 84199.664
$ws.Range("I42").Value = 87900
$ws.Range("J42").Value = 76799
$ws.Range("K42").Value = 87900
$ws.Range("L42").Value = 76799
$ws.Range("M42").Value = -87337
$ws.Range("N42").Value = -77925
# Row 49
$ws.Range("H49").Value = 84199.664
$ws.Range("I49").Value = 87900
$ws.Range("J49").Value = 76799
$ws.Range("K49").Value = 87900
$ws.Range("L49").Value = 76799
$ws.Range("M49").Value = -87753
$ws.Range("N49").Value = -77093
# Row 82
$ws.Range("H82").Value = 5016.4287
$ws.Range("I82").Value = 5823
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 5823
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -5462
$ws.Range("N82").Value = -3722
# Row 85
$ws.Range("H85").Value = 5016.4287
$ws.Range("I85").Value = 5823
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 5823
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -4575
$ws.Range("N85").Value = -5496
# Row 100
$ws.Range("H100").Value = 5784.857
$ws.Range("I100").Value = 2624.25
$ws.Range("K100").Value = 2624.25
$ws.Range("M100").Value = -2083.25
# Row 122
$ws.Range("H122").Value = 4779.0586
$ws.Range("I122").Value = 5408.154
$ws.Range("K122").Value = 16224.462
$ws.Range("M122").Value = -13774.462
# Row 132
$ws.Range("H132").Value = 535912.2
$ws.Range("I132").Value = 748252.4399999999
$ws.Range("K132").Value = 2244757.32
$ws.Range("M132").Value = -2242227.32

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2849.8823
$ws.Range("I122").Value = 1538.7693
$ws.Range("J122").Value = 7111
$ws.Range("K122").Value = 4616.3079
$ws.Range("L122").Value = 21333
$ws.Range("M122").Value = -2166.3079
$ws.Range("N122").Value = -26233
# Row 126
$ws.Range("H126").Value = 31789.715
$ws.Range("J126").Value = 6096.75
$ws.Range("L126").Value = 18290.25
$ws.Range("N126").Value = -23230.25
# Row 132
$ws.Range("H132").Value = 9532.5
$ws.Range("J132").Value = 4399
$ws.Range("L132").Value = 13197
$ws.Range("N132").Value = -18257
